# Auto-generated script applying the "想去人数"/"最低票价" updates described in the diff.
# Sheet order in the workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wb = $excel.ActiveWorkbook

# --- Worksheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1469
$ws.Range("F4").Value = 2111
$ws.Range("F5").Value = 7226
$ws.Range("F6").Value = 580
$ws.Range("G6").Value = "不可售"
$ws.Range("F7").Value = 64
$ws.Range("F8").Value = 4700
$ws.Range("F9").Value = 6929
$ws.Range("F11").Value = 257
$ws.Range("F12").Value = 1455
$ws.Range("F13").Value = 843
$ws.Range("F14").Value = 141
$ws.Range("F15").Value = 41
$ws.Range("F16").Value = 1157
$ws.Range("F18").Value = 145
$ws.Range("F20").Value = 207
$ws.Range("F22").Value = 1121
$ws.Range("F24").Value = 44
$ws.Range("F25").Value = 1200
$ws.Range("F26").Value = 38
$ws.Range("F29").Value = 39
$ws.Range("F30").Value = 145
$ws.Range("F31").Value = 12
$ws.Range("F32").Value = 31
$ws.Range("F33").Value = 61
$ws.Range("F34").Value = 25
$ws.Range("F36").Value = 540
$ws.Range("F37").Value = 418
$ws.Range("F40").Value = 346
$ws.Range("F42").Value = 563
$ws.Range("F43").Value = 76
$ws.Range("F44").Value = 134
$ws.Range("F46").Value = 16

# --- Worksheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F17").Value = 554
$ws.Range("F18").Value = 17
$ws.Range("F19").Value = 8
$ws.Range("F28").Value = 20
$ws.Range("F32").Value = 974
$ws.Range("F33").Value = 600
$ws.Range("F36").Value = 103
$ws.Range("F44").Value = 5

# --- Worksheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 722
$ws.Range("F6").Value = 649
$ws.Range("F8").Value = 1492
$ws.Range("F9").Value = 2345

# --- Worksheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 722
$ws.Range("F4").Value = 1469
$ws.Range("F7").Value = 649
$ws.Range("F8").Value = 649
$ws.Range("F10").Value = 7227
$ws.Range("F11").Value = 64
$ws.Range("F12").Value = 4700
$ws.Range("F14").Value = 6929
$ws.Range("F15").Value = 257
$ws.Range("F16").Value = 1455
$ws.Range("F18").Value = 554
$ws.Range("F19").Value = 843
$ws.Range("F20").Value = 41
$ws.Range("F21").Value = 1157
$ws.Range("F22").Value = 145
$ws.Range("F23").Value = 1121
$ws.Range("F26").Value = 44
$ws.Range("F27").Value = 1200
$ws.Range("F28").Value = 20
$ws.Range("F31").Value = 31
$ws.Range("F32").Value = 61
$ws.Range("F33").Value = 25
$ws.Range("F34").Value = 974
$ws.Range("F35").Value = 540
$ws.Range("F36").Value = 600
$ws.Range("F37").Value = 418
$ws.Range("F40").Value = 103
$ws.Range("F41").Value = 346
$ws.Range("F42").Value = 563
$ws.Range("F46").Value = 134
$ws.Range("F49").Value = 16

